$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of column R (2021) into the new column S (2022)
# so the new cells pick up identical number formats / fonts / borders.
$ws.Range("R3:R8").Copy()
$ws.Range("S3:S8").PasteSpecial(-4122)

# New "2022" column header
$ws.Range("S3").Value = 2022

# Row 4: "a) Number of branches ... per 100 000 adults"
# Formulas (R6/R8*100000) replaced by static values, and a new 2022 value added.
$ws.Range("R4").Value = 6.9132648934880807
$ws.Range("S4").Value = 6.9031689452913012

# Row 5: "b) Number of ATMs ... per 100 000 adults"
# Formulas (R7/R8*100000) replaced by static values, and a new 2022 value added.
$ws.Range("R5").Value = 42.321589572314856
$ws.Range("S5").Value = 44.306188104841333

# Row 6: Total branches of commercial banks
$ws.Range("S6").Value = 318

# Row 7: Total ATMs
$ws.Range("S7").Value = 2041

# Row 8: Adult resident population
$ws.Range("R8").Value = 4513063
$ws.Range("S8").Value = 4606580

# Restore the originally selected cell (shifted from R15 to R13)
$ws.Range("R13").Select()
